$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125, pushing the existing row 125 (and below) down to 126.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new price record.
$ws.Cells.Item(125, 1).Value = 5
$ws.Cells.Item(125, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(125, 3).Value = "Maule"
$ws.Cells.Item(125, 4).Value = 44628
$ws.Cells.Item(125, 5).Value = 7
$ws.Cells.Item(125, 6).Value = 100112030
$ws.Cells.Item(125, 7).Value = "Poroto granado"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 400
$ws.Cells.Item(125, 11).Value = 23000
$ws.Cells.Item(125, 12).Value = 23000
$ws.Cells.Item(125, 13).Value = 23000
$ws.Cells.Item(125, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(125, 15).Value = "Región del Maule"
$ws.Cells.Item(125, 16).Value = 920
$ws.Cells.Item(125, 17).Value = 25
$ws.Cells.Item(125, 18).Value = "Hortaliza"
